# Add 2022-Q3 data
# 1) Insert a new worksheet named "2022-Q3" right after "总计" (so it becomes the 2nd tab).
# 2) Populate it with the quarterly fund-holdings table.
# 3) Insert a new summary row in "总计" for 2022-Q3, shifting the existing rows down.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- Create the new "2022-Q3" sheet right after "总计" ---
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Copy the cell formatting (header band + bordered index column) from the existing
# "2022-Q2" sheet so the new sheet matches the workbook's established look.
$q2Sheet.Range("A1:H8").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

# --- Header row ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Data rows ---
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'515210"
$newSheet.Range("C2").Value = "国泰中证钢铁ETF"
$newSheet.Range("D2").Value = "'14.23"
$newSheet.Range("E2").Value = "'97.88"
$newSheet.Range("F2").Value = "'2.86"
$newSheet.Range("G2").Value = "'0.4070"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'502023"
$newSheet.Range("C3").Value = "鹏华国证钢铁行业指数（LOF）A"
$newSheet.Range("D3").Value = "'9.48"
$newSheet.Range("E3").Value = "'94.49"
$newSheet.Range("F3").Value = "'2.57"
$newSheet.Range("G3").Value = "'0.2436"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'012810"
$newSheet.Range("C4").Value = "鹏华国证钢铁行业指数（LOF）C"
$newSheet.Range("D4").Value = "'4.34"
$newSheet.Range("E4").Value = "'94.49"
$newSheet.Range("F4").Value = "'2.57"
$newSheet.Range("G4").Value = "'0.1115"
$newSheet.Range("H4").Value = 10

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'168203"
$newSheet.Range("C5").Value = "中融国证钢铁行业指数A"
$newSheet.Range("D5").Value = "'3.34"
$newSheet.Range("E5").Value = "'92.81"
$newSheet.Range("F5").Value = "'2.53"
$newSheet.Range("G5").Value = "'0.0845"
$newSheet.Range("H5").Value = 10

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'013802"
$newSheet.Range("C6").Value = "财通资管中证钢铁指数A"
$newSheet.Range("D6").Value = "'0.08"
$newSheet.Range("E6").Value = "'92.45"
$newSheet.Range("F6").Value = "'2.80"
$newSheet.Range("G6").Value = "'0.0022"
$newSheet.Range("H6").Value = 9

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'013803"
$newSheet.Range("C7").Value = "财通资管中证钢铁指数C"
$newSheet.Range("D7").Value = "'0.01"
$newSheet.Range("E7").Value = "'92.45"
$newSheet.Range("F7").Value = "'2.80"
$newSheet.Range("G7").Value = "'0.0003"
$newSheet.Range("H7").Value = 9

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'016815"
$newSheet.Range("C8").Value = "中融国证钢铁行业指数C"
$newSheet.Range("D8").Value = "'0.00"
$newSheet.Range("E8").Value = "'92.81"
$newSheet.Range("F8").Value = "'2.53"
$newSheet.Range("G8").Value = 0
$newSheet.Range("H8").Value = 10

# --- Update "总计" summary sheet: push existing data rows down one and insert the
#     new 2022-Q3 row at the top of the data (row 2). Done by copying bottom-up so
#     nothing is overwritten before it's read. ---
for ($r = 7; $r -ge 2; $r--) {
    $destRow = $r + 1
    $totalSheet.Cells.Item($destRow, 1).Value = $totalSheet.Cells.Item($r, 1).Value()
    $totalSheet.Cells.Item($destRow, 2).Value = $totalSheet.Cells.Item($r, 2).Value()
    $totalSheet.Cells.Item($destRow, 3).Value = $totalSheet.Cells.Item($r, 3).Value()
    $totalSheet.Cells.Item($destRow, 4).Value = $totalSheet.Cells.Item($r, 4).Value()
}

# Row 8 (now holding the old row-7 "2020-Q4" data) needs the bordered index-column
# style restored on A8, since it was a brand-new cell with no prior formatting.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A8").PasteSpecial(-4122)   # xlPasteFormats

# Write the new 2022-Q3 summary values into row 2.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.85
